# Applies the textual edits described by the commit "Fixed issue number 5"
# 1. Inserts a new explanatory sentence about binary vectors into section 2.1.
# 2. Adds a missing space before the "[5]" citation ("input[5]" -> "input [5]").
#
# NOTE: The remaining hunks in the source diff are pure Word pagination
# artifacts (w:lastRenderedPageBreak relocation / run re-splitting caused by
# the reflow after the text insertion above) and carry no visible text
# changes, so Word's own repagination after these edits reproduces them;
# they are not applied manually here.

$d = $word.ActiveDocument

# --- Edit 1: insert the new sentence about binary vectors -------------------
$findText = "some form of data. The SDR encoder of the HTM"
$replaceText = "some form of data. A binary vector refers to a vector with binary values (0 or 1), which typically has a fixed length, and each element of the vector corresponds to a feature or attribute of the input data. If a feature is present in the input data, its corresponding element in the binary vector is set to 1, otherwise it is set to 0. The SDR encoder of the HTM"

$d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null

# --- Edit 2: fix spacing of the "[5]" citation -------------------------------
$findText2 = "representation of the input[5]."
$replaceText2 = "representation of the input [5]."

$d.Content.Find.Execute($findText2, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText2, 2) | Out-Null

Write-Host "Edits applied"
